$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.922.91"
$ws.Range("E2").Value = "  +0.76%  "

$ws.Range("D3").Value = "1.758.03"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3841"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3384"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07204"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.134"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.151"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").Value = "1.756.91"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001056"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06612"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.61%  "

$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.213"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.60%  "

$ws.Range("D23").Value = "27.951.03"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.380"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.303"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.20%  "

$ws.Range("D29").Value = "1.960.26"
$ws.Range("E29").Value = "  -1.22%  "

$ws.Range("E30").Value = "  -15.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.018"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.820"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08817"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6597"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06190"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02291"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.151"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.33%  "

$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.502"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.42%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2105"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.207"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.947"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.19%  "

$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6021"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.93%  "

$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.122"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.170"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "

